$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Paragraph "return the good number" — drop the stray _GoBack bookmark
#    that used to sit between the two runs (it has moved further down, see
#    step 3 below).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(27)
if ($p1.Range.Text -notmatch "return the good number") {
    throw "Unexpected paragraph 27 text: [$($p1.Range.Text)]"
}
$xml1 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="9"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t>return the good number</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
"@
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Empty paragraph right after "GenerateRandomDieValue()..." becomes a new
#    bulleted requirement: "use random object (rnd) as required"
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(34)
if ($p2.Range.Text.TrimEnd([char]13) -ne "") {
    throw "Unexpected paragraph 34 text: [$($p2.Range.Text)]"
}
$xml2 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="9"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t>use random object (</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t>rnd</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t>) as required</w:t>
  </w:r>
</w:p>
"@
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) First of the three centred blank paragraphs before "SCREEN SAMPLES"
#    becomes a new bulleted requirement: "use while loop" (carrying the
#    relocated _GoBack bookmark) and a fresh blank list paragraph is added
#    right after it.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(36)
if ($p3.Range.Text.TrimEnd([char]13) -ne "") {
    throw "Unexpected paragraph 36 text: [$($p3.Range.Text)]"
}
$xml3 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="9"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">use </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">while </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
    <w:t>loop</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$p3.Range.InsertXML($xml3)

Write-Output "done: $($d.Paragraphs.Count) paragraphs"
